# Adding the changes we made on may 9th
# Appends rows 22-31 of walkingToRunning sensor data to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2000, "walkingToRunning", 25.90989303588867, 1.824093818664551, -40.16428375244141, 0.1055892638440206, 4.193923464361182, 1.237120070547445),
    @(2100, "walkingToRunning", -36.62774276733398, -9.502298355102541, -32.28135681152344, -12.19494866425143, 0.3453140144921644, 11.02819687465457),
    @(2200, "walkingToRunning", -5.977139472961426, -25.55127906799316, -0.151648998260498, 3.551830255760532, -8.829065172880417, 8.992697445851405),
    @(2300, "walkingToRunning", -3.294158458709717, 12.39510154724121, -7.373503684997559, 4.123927570738847, -6.680597377273238, 5.334880648918858),
    @(2400, "walkingToRunning", 48.08645629882812, -14.49578666687012, 15.91215991973877, 4.32835804291493, 5.768243861648356, 1.005123424080155),
    @(2500, "walkingToRunning", -71.25363159179688, -23.40527725219727, -48.95271301269531, -0.7777560342033203, 4.314871306689256, 2.766960350972315),
    @(2600, "walkingToRunning", -2.450338363647461, 4.339826583862305, -12.33499145507812, -5.723839211014095, 1.048127250851321, 5.411107558124477),
    @(2700, "walkingToRunning", 4.199068546295166, -20.79486465454102, -6.71635627746582, -2.712708458585516, 9.354416433370432, 0.1980249521866684),
    @(2800, "walkingToRunning", 16.30048179626465, 1.202856779098511, 33.71396636962891, 1.900770921752156, -1.904725146743598, -0.06296492072769189),
    @(2900, "walkingToRunning", 16.76902770996094, 16.37147903442383, 8.113217353820801, 6.883454844636765, -7.698698745583567, -1.359904374716438)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

Write-Output "appended $($data.Length) rows starting at row $startRow"
